$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Bowling Green Ohio" (style index 1)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bowling Green Ohio")

$ws.Range("E2").Value = 0.0784
$ws.Range("E3").Value = 0.0784
$ws.Range("E4").Value = 0.0784

$ws.Range("L4").Value = 0.0196
$ws.Range("M4").Value = 0.02
$ws.Range("N4").Value = 0.0392
$ws.Range("O4").Value = 0.0417
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws.Range("O7").ClearContents()

$ws.Range("K10").Value = 0.0115
$ws.Range("L10").Value = 0.0119
$ws.Range("N10").Value = 0.0232
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0.003325
$ws.Range("Q10").Value = 0.003325
$ws.Range("R10").Value = 0.009975
$ws.Range("S10").Value = 0.003325
$ws.Range("T10").Value = 0.003325
$ws.Range("U10").Value = 0.003325
$ws.Range("V10").Value = 0.009975
$ws.Range("W10").Value = 0.0399

# ---------------------------------------------------------------------------
# Sheet "Hyderabad India" (style index 2)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hyderabad India")

$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

# ---------------------------------------------------------------------------
# Sheet "Langley Canada" (style index 4)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Langley Canada")

$ws.Range("E2").Value = 0.0233
$ws.Range("E3").Value = 0.0233
$ws.Range("E4").Value = 0.0233

$ws.Range("O4").Value = 0.0119
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws.Range("O7").ClearContents()

$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0.00588333333333333
$ws.Range("Q10").Value = 0.00588333333333333
$ws.Range("R10").Value = 0.01765
$ws.Range("S10").Value = 0.00588333333333333
$ws.Range("T10").Value = 0.00588333333333333
$ws.Range("U10").Value = 0.00588333333333333
$ws.Range("V10").Value = 0.01765
$ws.Range("W10").Value = 0.0706

# ---------------------------------------------------------------------------
# Sheet "Las Vegas Nevada" (style index 5)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Las Vegas Nevada")

$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

$ws.Range("O7").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Apodaca Pmc Plant 2 Mexico" (style index 8)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Apodaca Pmc Plant 2 Mexico")

$ws.Range("E2").Value = 0.1852
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0

$ws.Range("E3").Value = 0.1896
$ws.Range("K3").Value = 0.1429
$ws.Range("M3").Value = 0.0476
$ws.Range("N3").Value = 0.1905
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.0270833333333333
$ws.Range("Q3").Value = 0.0270833333333333
$ws.Range("R3").Value = 0.08125
$ws.Range("S3").Value = 0.0270833333333333
$ws.Range("T3").Value = 0.0270833333333333
$ws.Range("U3").Value = 0.0270833333333333
$ws.Range("V3").Value = 0.08125
$ws.Range("W3").Value = 0.325
